$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ==========================================================================
# Add new attendance columns F..Y (per-day marks, plus Present/Total/%)
# to the existing 3-row / 5-column sheet, growing it to A1:Y3.
# ==========================================================================

# --- Header row (row 1): new column labels F1..Y1 ------------------------
# Note: the "YYYY-MM-DD" looking header labels (H1:V1) must stay literal
# text, exactly like the rest of the sheet, not get auto-converted into
# date serial numbers by Excel's smart-entry parsing. Prefixing with a
# leading apostrophe forces text entry; the apostrophe itself is not
# stored as part of the cell's value/text.
$headerValues = [ordered]@{
    "F1" = "2025-06-22_x.1";
    "G1" = "2025-06-22_y.1";
    "H1" = "'2025-06-08";
    "I1" = "'2025-06-09";
    "J1" = "'2025-06-10";
    "K1" = "'2025-06-11";
    "L1" = "'2025-06-12";
    "M1" = "'2025-06-13";
    "N1" = "'2025-06-14";
    "O1" = "'2025-06-15";
    "P1" = "'2025-06-16";
    "Q1" = "'2025-06-17";
    "R1" = "'2025-06-18";
    "S1" = "'2025-06-19";
    "T1" = "'2025-06-20";
    "U1" = "'2025-06-21";
    "V1" = "'2025-06-22";
    "W1" = "Present";
    "X1" = "Total";
    "Y1" = "Attendance %";
}

foreach ($addr in $headerValues.Keys) {
    $ws.Range($addr).Value = $headerValues[$addr]
}

# Apply the same look as the existing header cells (A1:E1) to the new
# header cells in one shot: bold font, thin box border, centered / top
# aligned text.
$headerRange = $ws.Range("F1:Y1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# --- Row 2 (Abhishek Pathak) daily marks + summary ------------------------
$row2Values = [ordered]@{
    "F2" = "✅"; "G2" = "✅"; "H2" = "❌"; "I2" = "❌"; "J2" = "❌";
    "K2" = "❌"; "L2" = "❌"; "M2" = "❌"; "N2" = "❌"; "O2" = "❌";
    "P2" = "❌"; "Q2" = "❌"; "R2" = "❌"; "S2" = "❌"; "T2" = "❌";
    "U2" = "❌"; "V2" = "✅";
}
foreach ($addr in $row2Values.Keys) {
    $ws.Range($addr).Value = $row2Values[$addr]
}
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 19
$ws.Range("Y2").Value = 15.8

# --- Row 3 (Shubham Pitekar) daily marks + summary ------------------------
$row3Values = [ordered]@{
    "F3" = "✅"; "G3" = "✅"; "H3" = "❌"; "I3" = "❌"; "J3" = "❌";
    "K3" = "❌"; "L3" = "❌"; "M3" = "❌"; "N3" = "❌"; "O3" = "❌";
    "P3" = "❌"; "Q3" = "❌"; "R3" = "❌"; "S3" = "❌"; "T3" = "❌";
    "U3" = "❌"; "V3" = "✅";
}
foreach ($addr in $row3Values.Keys) {
    $ws.Range($addr).Value = $row3Values[$addr]
}
$ws.Range("W3").Value = 5
$ws.Range("X3").Value = 19
$ws.Range("Y3").Value = 26.3
